# "Running suites a and b"
# Switch the Runmode column (C) on the "Test Cases" sheet from "N" to "Y"
# for rows 3-16, and leave the active selection on a single cell (C14)
# instead of the prior C3:C16 range selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("C3:C16").Value = "Y"

$ws.Range("C14").Select() | Out-Null
